$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# like "1.000" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.695.94"
$ws.Range("E2").Value = "  +3.02%  "
$ws.Range("D3").Value = "1.695.73"
$ws.Range("E3").Value = "  +1.68%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "317.34"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.3969"
$ws.Range("E7").Value = "  +2.61%  "
$ws.Range("D8").Value = "0.4083"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").Value = "1.505"
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("D10").Value = "1.000"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").Value = "50.96"
$ws.Range("E11").Value = "  -5.95%  "
$ws.Range("D12").Value = "0.08981"
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("D13").Value = "7.222"
$ws.Range("E13").Value = "  +5.62%  "
$ws.Range("D14").Value = "23.50"
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "8.198"
$ws.Range("E15").Value = "  +12.22%  "
$ws.Range("D16").Value = "0.00001332"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("D17").Value = "1.697.71"
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "100.18"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").Value = "0.07012"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "19.77"
$ws.Range("E20").Value = "  +4.02%  "
$ws.Range("D21").Value = "7.063"
$ws.Range("E21").Value = "  +7.98%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "14.24"
$ws.Range("E23").Value = "  +2.65%  "
$ws.Range("D24").Value = "24.691.43"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("D25").Value = "3.125"
$ws.Range("E25").Value = "  +5.96%  "
$ws.Range("D26").Value = "2.345"
$ws.Range("E26").Value = "  +1.12%  "
$ws.Range("D27").Value = "22.77"
$ws.Range("E27").Value = "  +4.53%  "
$ws.Range("D28").Value = "163.00"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "137.21"
$ws.Range("E29").Value = "  +4.94%  "
$ws.Range("D30").Value = "5.183"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Value = "7.443"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").Value = "1.881.58"
$ws.Range("E32").Value = "  +1.84%  "
$ws.Range("D33").Value = "1.084"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").Value = "0.08617"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "7.148"
$ws.Range("E35").Value = "  +1.04%  "
$ws.Range("D36").Value = "11.53"
$ws.Range("E36").Value = "  +5.59%  "
$ws.Range("D37").Value = "0.2754"
$ws.Range("E37").Value = "  +4.67%  "
$ws.Range("D38").Value = "1.924"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "14.48"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "0.09244"
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("D41").Value = "0.02732"
$ws.Range("E41").Value = "  +8.32%  "
$ws.Range("D42").Value = "1.481"
$ws.Range("E42").Value = "  +1.30%  "
$ws.Range("D43").Value = "0.7715"
$ws.Range("E43").Value = "  +3.47%  "
$ws.Range("D44").Value = "2.639"
$ws.Range("E44").Value = "  +9.86%  "
$ws.Range("D45").Value = "15.81"
$ws.Range("E45").Value = "  +5.86%  "
$ws.Range("D46").Value = "0.7201"
$ws.Range("E46").Value = "  +2.50%  "
$ws.Range("D47").Value = "4.229"
$ws.Range("E47").Value = "  +3.44%  "
$ws.Range("D48").Value = "1.000"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "141.11"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("B50").Value = "Flow"
$ws.Range("C50").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D50").Value = "1.330"
$ws.Range("E50").Value = "  +7.90%  "
$ws.Range("D51").Value = "0.07994"
$ws.Range("E51").Value = "  +2.06%  "
